# Apply updated transition-probability values to Sheet1 (Utah St. B team-specific matrix).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1825095057034221
$ws.Range("C2").Value = 0.5741444866920152
$ws.Range("P2").Value = 0.1482889733840304
$ws.Range("S2").Value = 0.09505703422053231

# Row 3
$ws.Range("B3").Value = 0.01290322580645161
$ws.Range("C3").Value = 0.03225806451612903
$ws.Range("J3").Value = 0.02580645161290323
$ws.Range("P3").Value = 0.7161290322580646
$ws.Range("S3").Value = 0.2129032258064516

# Row 4
$ws.Range("J4").Value = 0.02
$ws.Range("P4").Value = 0.62
$ws.Range("S4").Value = 0.36

# Row 6
$ws.Range("B6").Value = 0.06511627906976744
$ws.Range("D6").Value = 0.03255813953488372
$ws.Range("F6").Value = 0.07906976744186046
$ws.Range("J6").Value = 0.1906976744186047
$ws.Range("O6").Value = 0.004651162790697674
$ws.Range("Q6").Value = 0.1627906976744186
$ws.Range("R6").Value = 0.06046511627906977
$ws.Range("S6").Value = 0.4046511627906977

# Row 7
$ws.Range("B7").Value = 0.09216589861751152
$ws.Range("D7").Value = 0.02764976958525346
$ws.Range("E7").Value = 0.004608294930875576
$ws.Range("F7").Value = 0.05990783410138249
$ws.Range("J7").Value = 0.1428571428571428
$ws.Range("O7").Value = 0.009216589861751152
$ws.Range("Q7").Value = 0.2304147465437788
$ws.Range("R7").Value = 0.05529953917050692
$ws.Range("S7").Value = 0.3778801843317972

# Row 8
$ws.Range("B8").Value = 0.09823182711198428
$ws.Range("D8").Value = 0.01571709233791748
$ws.Range("F8").Value = 0.06483300589390963
$ws.Range("J8").Value = 0.0962671905697446
$ws.Range("O8").Value = 0.01768172888015717
$ws.Range("Q8").Value = 0.1944990176817289
$ws.Range("R8").Value = 0.068762278978389
$ws.Range("S8").Value = 0.444007858546169

# Row 9
$ws.Range("B9").Value = 0.1038961038961039
$ws.Range("D9").Value = 0.01948051948051948
$ws.Range("F9").Value = 0.06493506493506493
$ws.Range("J9").Value = 0.09090909090909091
$ws.Range("O9").Value = 0.03246753246753246
$ws.Range("Q9").Value = 0.2402597402597403
$ws.Range("R9").Value = 0.06493506493506493
$ws.Range("S9").Value = 0.3831168831168831

# Row 10
$ws.Range("B10").Value = 0.09671532846715329
$ws.Range("D10").Value = 0.02372262773722628
$ws.Range("F10").Value = 0.07846715328467153
$ws.Range("J10").Value = 0.1085766423357664
$ws.Range("O10").Value = 0.01551094890510949
$ws.Range("Q10").Value = 0.208029197080292
$ws.Range("R10").Value = 0.07755474452554745
$ws.Range("S10").Value = 0.3914233576642336

# Row 11
$ws.Range("G11").Value = 0.1242424242424242
$ws.Range("J11").Value = 0.1060606060606061
$ws.Range("K11").Value = 0.1818181818181818
$ws.Range("L11").Value = 0.5787878787878787
$ws.Range("S11").Value = 0.00909090909090909

# Row 12
$ws.Range("G12").Value = 0.7085427135678392
$ws.Range("J12").Value = 0.2060301507537688
$ws.Range("K12").Value = 0.005025125628140704
$ws.Range("L12").Value = 0.03015075376884422
$ws.Range("S12").Value = 0.05025125628140704

# Row 13
$ws.Range("G13").Value = 0.7241379310344828
$ws.Range("J13").Value = 0.2586206896551724
$ws.Range("S13").Value = 0.01724137931034483

# Row 15
$ws.Range("F15").Value = 0.02590673575129534
$ws.Range("H15").Value = 0.1813471502590674
$ws.Range("I15").Value = 0.05699481865284974
$ws.Range("J15").Value = 0.2901554404145077
$ws.Range("K15").Value = 0.1191709844559585
$ws.Range("M15").Value = 0.02072538860103627
$ws.Range("O15").Value = 0.0310880829015544
$ws.Range("S15").Value = 0.2746113989637305

# Row 16
$ws.Range("F16").Value = 0.02857142857142857
$ws.Range("H16").Value = 0.1942857142857143
$ws.Range("I16").Value = 0.06857142857142857
$ws.Range("J16").Value = 0.3371428571428571
$ws.Range("K16").Value = 0.1657142857142857
$ws.Range("M16").Value = 0.01714285714285714
$ws.Range("N16").Value = 0.005714285714285714
$ws.Range("O16").Value = 0.06857142857142857
$ws.Range("S16").Value = 0.1142857142857143

# Row 17
$ws.Range("F17").Value = 0.008968609865470852
$ws.Range("H17").Value = 0.2152466367713005
$ws.Range("I17").Value = 0.07174887892376682
$ws.Range("J17").Value = 0.3968609865470852
$ws.Range("K17").Value = 0.1098654708520179
$ws.Range("M17").Value = 0.0336322869955157
$ws.Range("N17").Value = 0.004484304932735426
$ws.Range("O17").Value = 0.05381165919282511
$ws.Range("S17").Value = 0.1053811659192825

# Row 18
$ws.Range("F18").Value = 0.006535947712418301
$ws.Range("H18").Value = 0.1699346405228758
$ws.Range("I18").Value = 0.09803921568627451
$ws.Range("J18").Value = 0.3660130718954248
$ws.Range("K18").Value = 0.1372549019607843
$ws.Range("M18").Value = 0.0392156862745098
$ws.Range("O18").Value = 0.0718954248366013
$ws.Range("S18").Value = 0.1111111111111111

# Row 19
$ws.Range("F19").Value = 0.01470588235294118
$ws.Range("H19").Value = 0.2638888888888889
$ws.Range("I19").Value = 0.06781045751633986
$ws.Range("J19").Value = 0.3439542483660131
$ws.Range("K19").Value = 0.113562091503268
$ws.Range("M19").Value = 0.02532679738562092
$ws.Range("O19").Value = 0.06372549019607843
$ws.Range("S19").Value = 0.1070261437908497
